$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures (text values, apostrophe-prefixed
# to preserve their original text/inline-string representation rather than being
# auto-converted to numeric/percentage values by Excel).
$ws.Range("D2").Value = "'320.71"
$ws.Range("E2").Value = "'-1.69%"
$ws.Range("D3").Value = "'39.30"
$ws.Range("E3").Value = "'-1.23%"
$ws.Range("E4").Value = "'12.42%"
$ws.Range("D5").Value = "'0.08006"
$ws.Range("E5").Value = "'-1.16%"
$ws.Range("E6").Value = "'-0.11%"
$ws.Range("D7").Value = "'1.903"
$ws.Range("E7").Value = "'-0.99%"
$ws.Range("D8").Value = "'0.9347"
$ws.Range("E8").Value = "'0.10%"
$ws.Range("D9").Value = "'0.1249"
$ws.Range("E9").Value = "'-4.58%"
$ws.Range("D10").Value = "'0.1951"
$ws.Range("E10").Value = "'-0.17%"
$ws.Range("D11").Value = "'8.737"
$ws.Range("E11").Value = "'27.32%"
$ws.Range("D12").Value = "'0.09125"
$ws.Range("E12").Value = "'-1.49%"
$ws.Range("D13").Value = "'0.03516"
$ws.Range("E13").Value = "'3.00%"
$ws.Range("D14").Value = "'0.09561"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("D15").Value = "'0.001293"
$ws.Range("E15").Value = "'-7.44%"
$ws.Range("D16").Value = "'0.006220"
$ws.Range("E16").Value = "'-2.94%"
$ws.Range("D17").Value = "'3.354"
$ws.Range("E17").Value = "'-0.11%"
$ws.Range("D18").Value = "'4.565"
$ws.Range("E18").Value = "'0.68%"
$ws.Range("D19").Value = "'2.950"
$ws.Range("E19").Value = "'-0.30%"
$ws.Range("D20").Value = "'0.3535"
$ws.Range("E20").Value = "'-0.01%"
$ws.Range("D21").Value = "'0.1429"
$ws.Range("E21").Value = "'6.53%"
$ws.Range("E22").Value = "'4.40%"
$ws.Range("D23").Value = "'0.04451"
$ws.Range("E23").Value = "'0.32%"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("E24").Value = "'3.39%"
$ws.Range("D25").Value = "'0.004369"
$ws.Range("E25").Value = "'0.28%"
$ws.Range("D26").Value = "'0.0001141"
$ws.Range("E26").Value = "'-11.53%"
$ws.Range("E27").Value = "'0.09%"
$ws.Range("D39").Value = "'0.02393"
$ws.Range("E39").Value = "'-2.90%"
$ws.Range("D40").Value = "'0.05169"
$ws.Range("E40").Value = "'-1.02%"
$ws.Range("D41").Value = "'0.007463"
$ws.Range("E41").Value = "'-3.04%"
$ws.Range("D42").Value = "'0.009198"
$ws.Range("E42").Value = "'7.77%"
$ws.Range("E43").Value = "'-2.03%"
$ws.Range("D44").Value = "'0.002122"
$ws.Range("E44").Value = "'0.58%"
$ws.Range("D45").Value = "'0.01115"
$ws.Range("E45").Value = "'37.16%"
$ws.Range("D46").Value = "'0.00006750"
$ws.Range("E46").Value = "'1.41%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'0.003007"
$ws.Range("E48").Value = "'5.48%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.06%"
